$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated team-specific transition-probability matrix values (Furman_B)
# Row 2
$ws.Range("B2").Value = 0.2321428571428572
$ws.Range("C2").Value = 0.4928571428571429
$ws.Range("J2").Value = 0.003571428571428571
$ws.Range("P2").Value = 0.1428571428571428
$ws.Range("S2").Value = 0.1285714285714286
# Row 3
$ws.Range("C3").Value = 0.02777777777777778
$ws.Range("J3").Value = 0.02083333333333333
$ws.Range("P3").Value = 0.7013888888888888
$ws.Range("S3").Value = 0.25
# Row 4
$ws.Range("J4").Value = 0.06976744186046512
$ws.Range("P4").Value = 0.6976744186046512
$ws.Range("S4").Value = 0.2325581395348837
# Row 6
$ws.Range("B6").Value = 0.07373271889400922
$ws.Range("D6").Value = 0.009216589861751152
$ws.Range("F6").Value = 0.05069124423963134
$ws.Range("J6").Value = 0.2672811059907834
$ws.Range("O6").Value = 0.009216589861751152
$ws.Range("Q6").Value = 0.1751152073732719
$ws.Range("R6").Value = 0.07834101382488479
$ws.Range("S6").Value = 0.336405529953917
# Row 7
$ws.Range("B7").Value = 0.07655502392344497
$ws.Range("D7").Value = 0.01435406698564593
$ws.Range("E7").Value = 0.004784688995215311
$ws.Range("F7").Value = 0.03827751196172249
$ws.Range("J7").Value = 0.1196172248803828
$ws.Range("O7").Value = 0.004784688995215311
$ws.Range("Q7").Value = 0.1818181818181818
$ws.Range("R7").Value = 0.1052631578947368
$ws.Range("S7").Value = 0.4545454545454545
# Row 8
$ws.Range("B8").Value = 0.08089887640449438
$ws.Range("D8").Value = 0.02022471910112359
$ws.Range("F8").Value = 0.06292134831460675
$ws.Range("J8").Value = 0.1033707865168539
$ws.Range("O8").Value = 0.008988764044943821
$ws.Range("Q8").Value = 0.197752808988764
$ws.Range("R8").Value = 0.09213483146067415
$ws.Range("S8").Value = 0.4337078651685393
# Row 9
$ws.Range("B9").Value = 0.09236947791164658
$ws.Range("D9").Value = 0.01606425702811245
$ws.Range("F9").Value = 0.04016064257028112
$ws.Range("J9").Value = 0.08835341365461848
$ws.Range("O9").Value = 0.02409638554216868
$ws.Range("Q9").Value = 0.1887550200803213
$ws.Range("R9").Value = 0.06425702811244979
$ws.Range("S9").Value = 0.4859437751004016
# Row 10
$ws.Range("B10").Value = 0.09355828220858896
$ws.Range("D10").Value = 0.01993865030674847
$ws.Range("E10").Value = 0.0007668711656441718
$ws.Range("F10").Value = 0.07438650306748466
$ws.Range("J10").Value = 0.1096625766871166
$ws.Range("O10").Value = 0.009202453987730062
$ws.Range("Q10").Value = 0.2523006134969325
$ws.Range("R10").Value = 0.0736196319018405
$ws.Range("S10").Value = 0.3665644171779141
# Row 11
$ws.Range("G11").Value = 0.1305732484076433
$ws.Range("J11").Value = 0.06369426751592357
$ws.Range("K11").Value = 0.2006369426751592
$ws.Range("L11").Value = 0.5828025477707006
$ws.Range("S11").Value = 0.02229299363057325
# Row 12
$ws.Range("G12").Value = 0.7121212121212122
$ws.Range("J12").Value = 0.1767676767676768
$ws.Range("L12").Value = 0.06565656565656566
$ws.Range("S12").Value = 0.04545454545454546
# Row 13
$ws.Range("G13").Value = 0.5490196078431373
$ws.Range("J13").Value = 0.392156862745098
$ws.Range("S13").Value = 0.05882352941176471
# Row 15
$ws.Range("F15").Value = 0.0310880829015544
$ws.Range("H15").Value = 0.1398963730569948
$ws.Range("I15").Value = 0.08290155440414508
$ws.Range("J15").Value = 0.4145077720207254
$ws.Range("K15").Value = 0.05181347150259067
$ws.Range("M15").Value = 0.03626943005181347
$ws.Range("O15").Value = 0.03626943005181347
$ws.Range("S15").Value = 0.2072538860103627
# Row 16
$ws.Range("F16").Value = 0.02409638554216868
$ws.Range("H16").Value = 0.2048192771084337
$ws.Range("I16").Value = 0.1144578313253012
$ws.Range("J16").Value = 0.3614457831325301
$ws.Range("K16").Value = 0.0963855421686747
$ws.Range("M16").Value = 0.01807228915662651
$ws.Range("N16").Value = 0.01807228915662651
$ws.Range("O16").Value = 0.04819277108433735
$ws.Range("S16").Value = 0.1144578313253012
# Row 17
$ws.Range("F17").Value = 0.01669758812615955
$ws.Range("H17").Value = 0.1725417439703154
$ws.Range("I17").Value = 0.1094619666048238
$ws.Range("J17").Value = 0.4174397031539889
$ws.Range("K17").Value = 0.09461966604823747
$ws.Range("M17").Value = 0.0111317254174397
$ws.Range("O17").Value = 0.04823747680890538
$ws.Range("S17").Value = 0.1298701298701299
# Row 18
$ws.Range("F18").Value = 0.01047120418848168
$ws.Range("H18").Value = 0.1989528795811518
$ws.Range("I18").Value = 0.08900523560209424
$ws.Range("J18").Value = 0.3717277486910995
$ws.Range("K18").Value = 0.1204188481675393
$ws.Range("M18").Value = 0.03141361256544502
$ws.Range("O18").Value = 0.05235602094240838
$ws.Range("S18").Value = 0.1256544502617801
# Row 19
$ws.Range("F19").Value = 0.01313868613138686
$ws.Range("H19").Value = 0.1832116788321168
$ws.Range("I19").Value = 0.1065693430656934
$ws.Range("J19").Value = 0.3715328467153285
$ws.Range("K19").Value = 0.1124087591240876
$ws.Range("M19").Value = 0.0218978102189781
$ws.Range("N19").Value = 0.00145985401459854
$ws.Range("O19").Value = 0.06350364963503649
$ws.Range("S19").Value = 0.1262773722627737
